$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "small_measure"
$ws.Range("A3").Value = "renovation"
$ws.Range("A4").Value = "demolition"
$ws.Range("A2:A4").Font.Bold = $false
$ws.Range("A5").Select() | Out-Null
